$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "datos actualizados" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 05:35"

# 2) Reorder country "Butan" in the list: it now appears right after "Comoras"
#    (row 205), pushing Montserrat / Seychelles / Groenlandia / Surinam down
#    one row each (rows 206-210), and update all stats for those rows.

# Row 145: Martinica - update Casos activos / Recuperados
$ws.Range("D145").Value = 91
$ws.Range("E145").Value = 82

# Row 178: Mongolia - update Casos activos / Recuperados
$ws.Range("D178").Value = 15
$ws.Range("E178").Value = 27

# Row 206: now Butan (was Montserrat)
$ws.Range("A206").Value = "Butan"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 2
$ws.Range("D206").Value = 5
$ws.Range("E206").Value = 6
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

# Row 207: now Montserrat (was Seychelles)
$ws.Range("A207").Value = "Montserrat"
$ws.Range("B207").Value = 11
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 8
$ws.Range("E207").Value = 2
$ws.Range("F207").Value = 1
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 1

# Row 208: now Seychelles (was Groenlandia)
$ws.Range("A208").Value = "Seychelles"
$ws.Range("B208").Value = 11
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 10
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209: now Groenlandia (was Surinam)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 210: now Surinam (was Butan)
$ws.Range("A210").Value = "Surinam"
$ws.Range("B210").Value = 10
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 9
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1
